$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 349.5
$ws.Range("I18").Value = 349.5
$ws.Range("K18").Value = 349.5
$ws.Range("M18").Value = -65.5

$ws.Range("H40").Value = 2200
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2200
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2200
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2550

$ws.Range("H42").Value = 796.7368
$ws.Range("I42").Value = 1316.5555
$ws.Range("J42").Value = 328.9
$ws.Range("K42").Value = 3949.6665
$ws.Range("L42").Value = 986.6999999999999
$ws.Range("M42").Value = -3719.6665
$ws.Range("N42").Value = -1446.7

$ws.Range("H70").Value = 1534.3334
$ws.Range("J70").Value = 1648.4667
$ws.Range("L70").Value = 4945.4001
$ws.Range("N70").Value = -5485.4001

$ws.Range("H73").Value = 1534.3334
$ws.Range("J73").Value = 1648.4667
$ws.Range("L73").Value = 4945.4001
$ws.Range("N73").Value = -6817.4001

$ws.Range("H133").Value = 34450
$ws.Range("J133").Value = 34450
$ws.Range("L133").Value = 34450
$ws.Range("N133").Value = -44570

$ws.Range("H134").Value = 50909.09
$ws.Range("J134").Value = 50909.09
$ws.Range("L134").Value = 50909.09
$ws.Range("N134").Value = -61049.09

$ws.Range("H137").Value = 1244.8
$ws.Range("I137").Value = 782.2857
$ws.Range("J137").Value = 1385.5652
$ws.Range("K137").Value = 2346.8571
$ws.Range("L137").Value = 4156.6956
$ws.Range("M137").Value = 203.1428999999998
$ws.Range("N137").Value = -9256.695599999999

$ws.Range("H138").Value = 7815149.5
$ws.Range("I138").Value = 3170.4443
$ws.Range("J138").Value = 10872011
$ws.Range("K138").Value = 9511.332900000001
$ws.Range("L138").Value = 32616033
$ws.Range("M138").Value = -4371.332900000001
$ws.Range("N138").Value = -32626313

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9407.857
$ws.Range("I2").Value = 12694
$ws.Range("K2").Value = 12694
$ws.Range("M2").Value = -12581

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H76").Value = 29500
$ws.Range("I76").Value = 9000
$ws.Range("K76").Value = 9000
$ws.Range("M76").Value = -8662

$ws.Range("H79").Value = 29500
$ws.Range("I79").Value = 9000
$ws.Range("K79").Value = 9000
$ws.Range("M79").Value = -7830

$ws.Range("H116").Value = 9407.857
$ws.Range("I116").Value = 12694
$ws.Range("K116").Value = 12694
$ws.Range("M116").Value = -10400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9407.857
$ws.Range("I3").Value = 12694
$ws.Range("K3").Value = 12694
$ws.Range("M3").Value = -12580

$ws.Range("H29").Value = 9766.666999999999
$ws.Range("I29").Value = 600
$ws.Range("J29").Value = 14350
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 14350
$ws.Range("M29").Value = -311
$ws.Range("N29").Value = -14928

$ws.Range("H80").Value = 791.2727
$ws.Range("I80").Value = 465.33334
$ws.Range("K80").Value = 465.33334
$ws.Range("M80").Value = 532.66666

$ws.Range("H83").Value = 791.2727
$ws.Range("I83").Value = 465.33334
$ws.Range("K83").Value = 2326.6667
$ws.Range("M83").Value = 2665.3333

$ws.Range("H140").Value = 61667
$ws.Range("J140").Value = 61667
$ws.Range("L140").Value = 61667
$ws.Range("N140").Value = -72027

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7525000.5
$ws.Range("I6").Value = 7525000.5
$ws.Range("K6").Value = 7525000.5
$ws.Range("M6").Value = -7524887.5

$ws.Range("H22").Value = 252.83333
$ws.Range("I22").Value = 252.83333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 252.83333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 97.16667000000001
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 2120.1082
$ws.Range("I31").Value = 1352.4117
$ws.Range("J31").Value = 2772.65
$ws.Range("K31").Value = 1352.4117
$ws.Range("L31").Value = 2772.65
$ws.Range("M31").Value = -1057.4117
$ws.Range("N31").Value = -3362.65

$ws.Range("H34").Value = 2120.1082
$ws.Range("I34").Value = 1352.4117
$ws.Range("J34").Value = 2772.65
$ws.Range("K34").Value = 1352.4117
$ws.Range("L34").Value = 2772.65
$ws.Range("M34").Value = -1150.4117
$ws.Range("N34").Value = -3176.65

$ws.Range("H57").Value = 9659.200000000001
$ws.Range("J57").Value = 9659.200000000001
$ws.Range("L57").Value = 9659.200000000001
$ws.Range("N57").Value = -10779.2

$ws.Range("H97").Value = 29599.5
$ws.Range("J97").Value = 29599.5
$ws.Range("L97").Value = 29599.5
$ws.Range("N97").Value = -31581.5

$ws.Range("H125").Value = 13386
$ws.Range("J125").Value = 13386
$ws.Range("L125").Value = 13386
$ws.Range("N125").Value = -18306

$ws.Range("H134").Value = 4705.1113
$ws.Range("I134").Value = 1222.5
$ws.Range("J134").Value = 7491.2
$ws.Range("K134").Value = 3667.5
$ws.Range("L134").Value = 22473.6
$ws.Range("M134").Value = -1132.5
$ws.Range("N134").Value = -27543.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1394.1
$ws.Range("I5").Value = 600
$ws.Range("J5").Value = 1658.8
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 4976.4
$ws.Range("M5").Value = -1688
$ws.Range("N5").Value = -5200.4

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws.Range("H131").Value = 1094.71
$ws.Range("I131").Value = 489
$ws.Range("J131").Value = 1147.3805
$ws.Range("K131").Value = 1467
$ws.Range("L131").Value = 3442.1415
$ws.Range("M131").Value = 3573
$ws.Range("N131").Value = -13522.1415

$ws.Range("H135").Value = 1394.1
$ws.Range("I135").Value = 600
$ws.Range("J135").Value = 1658.8
$ws.Range("K135").Value = 5400
$ws.Range("L135").Value = 14929.2
$ws.Range("M135").Value = -2865
$ws.Range("N135").Value = -19999.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2195.6
$ws.Range("I126").Value = 1382.4
$ws.Range("J126").Value = 2466.6667
$ws.Range("K126").Value = 4147.200000000001
$ws.Range("L126").Value = 7400.000100000001
$ws.Range("M126").Value = -1677.200000000001
$ws.Range("N126").Value = -12340.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1400
$ws.Range("I22").Value = 1400
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1400
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1105
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 1400
$ws.Range("I27").Value = 1400
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1400
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1293
$ws.Range("N27").ClearContents()

$ws.Range("H46").Value = 489.83334
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 489.83334
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 489.83334
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -865.83334

$ws.Range("H125").Value = 29753.666
$ws.Range("I125").Value = 650
$ws.Range("J125").Value = 38069
$ws.Range("K125").Value = 650
$ws.Range("L125").Value = 38069
$ws.Range("M125").Value = 4270
$ws.Range("N125").Value = -47909

$ws.Range("H132").Value = 5593.5835
$ws.Range("I132").Value = 4047.1667
$ws.Range("J132").Value = 7140
$ws.Range("K132").Value = 12141.5001
$ws.Range("L132").Value = 21420
$ws.Range("M132").Value = -9611.500100000001
$ws.Range("N132").Value = -26480

$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 251875.25
$ws.Range("I126").Value = 500750.5
$ws.Range("K126").Value = 1502251.5
$ws.Range("M126").Value = -1499781.5

$ws.Range("H132").Value = 19233740
$ws.Range("I132").Value = 33335912
$ws.Range("J132").Value = 3505
$ws.Range("K132").Value = 100007736
$ws.Range("L132").Value = 10515
$ws.Range("M132").Value = -100005206
$ws.Range("N132").Value = -15575

Write-Output "Applied Titan_Profits updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
